$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("1ST Q 2024")
$ws4 = $wb.Worksheets.Item("4th Q 2024")

$ws1.Activate()

# New travel entries for the 1st quarter sheet (rows 15-19).
# Shared-string insertion order matters, so the cells are written in the
# same order the author must have typed them (row15-A, row16-B, row15-B,
# row16-A, row17-A, row17-B, row18-A, row18-B, row19-A, row19-B).
$ws1.Range("A15").Value = "MALANAN, ALMA A."
$ws1.Range("B16").Value = "HRMO"
$ws1.Range("B15").Value = "CITY ADMINISTRATOR"
$ws1.Range("A16").Value = "SUSA, NANETTE B."
$ws1.Range("A17").Value = "BAYOT, JENNIFER"
$ws1.Range("B17").Value = "OIC- GSO"
$ws1.Range("A18").Value = "AMON, RHEALYN OCAMPO"
$ws1.Range("B18").Value = "ACCOUNTANT IV"
$ws1.Range("A19").Value = "LERIO, ROSEMARIE VERGARA"
$ws1.Range("B19").Value = "CITY ACCOUNTANT"

$ws1.Range("C15:C19").Value = "TAIWAN"
$ws1.Range("D15:D19").Value = "FEBRUARY 22- 26, 2024"
$ws1.Range("E15:E19").Value = "PERSONAL"

# Copy the "/" + "1" formatting/values down from row 14 so the new rows 15
# pick up the same styles as the rest of the table.
$ws1.Range("H14").Copy($ws1.Range("H15"))
$ws1.Range("I15").Value = 1

$ws1.Rows("15:19").RowHeight = 30

$ws1.Range("O10").Select()

$ws4.Range("E19").Select()

$ws1.Activate()
